$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @("0.1380735117988072", "0.189336789175735")
    3  = @("0.1346415758165991", "0.200680664247515")
    4  = @("0.04629597091307004", "0.6612275134895177")
    5  = @("0.1803593964616055", "0.08535381217266605")
    6  = @("0.1275116976770918", "0.2283991660541867")
    7  = @("-0.1778420293220559", "0.09169466198272562")
    8  = @("0.06879948303560469", "0.5146259231495384")
    9  = @("0.4554130660137468", "5.082960642748461e-06")
    10 = @("-0.320666183102976", "0.001830726413273231")
    11 = @("0.5387839372184192", "3.020070684012948e-08")
    12 = @("0.04512336699790357", "0.6710441167508682")
    13 = @("-0.2002122069313494", "0.0570647328831485")
    14 = @("-0.353373337670157", "0.0005497190329492098")
    15 = @("0.6267782526274787", "2.332215500968583e-11")
    16 = @("0.9836608241155468", "8.165814046459598e-69")
    17 = @("-0.5649690258760608", "5.420012969870346e-09")
    18 = @("0.3648624212277486", "0.0003767401393740044")
    19 = @("0.5827091103973073", "1.102073402945877e-09")
    20 = @("0.5505092224414716", "1.307276354452293e-08")
    21 = @("-0.6036968385343067", "2.384556196254475e-10")
    22 = @("0.4473763460813174", "8.749363260680295e-06")
    23 = @("0.9892887999456466", "5.155908730501901e-77")
    24 = @("-0.4243747729830099", "2.761998493444329e-05")
    25 = @("0.2757317777168685", "0.008158355348982125")
    26 = @("0.5231400147403523", "8.792810686124908e-08")
    27 = @("-0.4243209457586143", "2.769171175221538e-05")
    28 = @("-0.5060559429376476", "3.100961730853606e-07")
    29 = @("0.4097066303405117", "5.510397502678359e-05")
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = [double]$pair[0]
    $ws.Cells.Item($row, 4).Value = [double]$pair[1]
}
